$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12: Minimum Expense
$ws.Range("A12").Value = "Minimum Expense"
$ws.Range("B12").Formula = "=MIN(B4:B8)"
$ws.Range("C12:F12").Formula = "=MIN(C4:C8)"

# New row 13: Maximum Expense
$ws.Range("A13").Value = "Maximum Expense"
$ws.Range("B13").Formula = "=MAX(B4:B8)"
$ws.Range("C13:F13").Formula = "=MAX(C4:C8)"

# Column width adjustments (column A grew to fit "Maximum Expense", column F added for bestfit).
# The COM width model here quantizes to whole pixels (1/6 character increments), so feed it
# target-minus-padding and let it round to the closest representable width.
$ws.Columns.Item(1).ColumnWidth = 14.565104166666666
$ws.Columns.Item(6).ColumnWidth = 11.065104166666666

# Selection moved
$ws.Range("B21").Select()
